$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows right after the header row (row 1), pushing all
# existing song rows (old rows 2-190) down to rows 10-198. Inserting one
# row at a time at index 2 (instead of a single multi-row range insert)
# is what reliably shifts the whole block down here.
for ($i = 0; $i -lt 8; $i++) {
    $ws.Rows.Item(2).Insert()
}

# The inserted rows pick up the formatting of the row above (the header),
# but in the target workbook these rows are plain/unstyled data rows, so
# clear any inherited formatting (one row at a time - the "2:9" string
# range form does not reliably apply here).
for ($r = 2; $r -le 9; $r++) {
    $ws.Rows.Item($r).ClearFormats()
}

# Fill in the new songs that were added at the top of the list.
$ws.Cells.Item(2, 2).Value = "올인(All In) - MV_처음 그날처럼 (2003)"
$ws.Cells.Item(2, 3).Value = "https://www.youtube.com/watch?v=ddD9G7KQzx0"

$ws.Cells.Item(3, 2).Value = "Yalın - Yeniden"
$ws.Cells.Item(3, 3).Value = "https://www.youtube.com/watch?v=iGut_MVMcUY"

$ws.Cells.Item(4, 2).Value = "Rafet El Roman & Derya - Özledim (Düet)"
$ws.Cells.Item(4, 3).Value = "https://www.youtube.com/watch?v=JJ1fR1X4NYk"

$ws.Cells.Item(5, 2).Value = "Yalın - Zalim (Official Video)"
$ws.Cells.Item(5, 3).Value = "https://www.youtube.com/watch?v=kPM5VXhpCfA"

$ws.Cells.Item(6, 2).Value = "Rafet El Roman ft. Sinem - Seni Seviyorum"
$ws.Cells.Item(6, 3).Value = "https://www.youtube.com/watch?v=B3OcAOzFOCc"

$ws.Cells.Item(7, 2).Value = "Rafet El Roman - Senden Sonra"
$ws.Cells.Item(7, 3).Value = "https://www.youtube.com/watch?v=Z2g8NAg1bbI"

$ws.Cells.Item(8, 2).Value = "Rafet El Roman - Kalbine Sürgün Feat. Ezo"
$ws.Cells.Item(8, 3).Value = "https://www.youtube.com/watch?v=7I3h7czMJeI"

$ws.Cells.Item(9, 2).Value = "Rafet El Roman & Derya - Unuturum Elbet"
$ws.Cells.Item(9, 3).Value = "https://www.youtube.com/watch?v=ScZFzmN-8XY"

Write-Output "inserted 8 rows and populated new songs"
